$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.462.64"
$ws.Range("E2").Value = "  -2.74%  "

$ws.Range("D3").Value = "1.973.42"
$ws.Range("E3").Value = "  -3.94%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.50"
$ws.Range("E5").Value = "  +0.88%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.633"
$ws.Range("E6").Value = "  -4.53%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "56.76"
$ws.Range("E7").Value = "  +4.65%  "

$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "58.22"
$ws.Range("E9").Value = "  -0.04%  "

$ws.Range("E10").Value = "  -0.91%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0729"
$ws.Range("E11").Value = "  -2.75%  "

$ws.Range("E12").Value = "  -3.09%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.943"
$ws.Range("E13").Value = "  +3.89%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.36"
$ws.Range("E14").Value = "  -2.23%  "

$ws.Range("D15").Value = "2.258.38"
$ws.Range("E15").Value = "  -4.17%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.24"
$ws.Range("E16").Value = "  -2.68%  "

$ws.Range("D17").Value = "1.980.33"
$ws.Range("E17").Value = "  -3.73%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.61"
$ws.Range("E18").Value = "  +4.92%  "

$ws.Range("D19").Value = "35.374.38"
$ws.Range("E19").Value = "  -2.83%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.29"
$ws.Range("E20").Value = "  -0.78%  "

$ws.Range("E21").Value = "  -2.12%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "231.86"
$ws.Range("E22").Value = "  -2.76%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.12"
$ws.Range("E23").Value = "  -2.16%  "

$ws.Range("E24").Value = "  +0.06%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.51"
$ws.Range("E25").Value = "  +18.45%  "

$ws.Range("E26").Value = "  -2.70%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "164.15"
$ws.Range("E27").Value = "  +0.08%  "

$ws.Range("E28").Value = "  -3.40%  "

$ws.Range("E29").Value = "  -5.01%  "

$ws.Range("E30").Value = "  -3.41%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.84"
$ws.Range("E31").Value = "  -3.96%  "

$ws.Range("E32").Value = "  -8.59%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0942"
$ws.Range("E33").Value = "  +14.95%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0586"
$ws.Range("E34").Value = "  -0.90%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.29"
$ws.Range("E35").Value = "  -3.67%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.36"
$ws.Range("E36").Value = "  +7.96%  "

$ws.Range("E37").Value = "  -0.07%  "

$ws.Range("E38").Value = "  -4.04%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.16"
$ws.Range("E39").Value = "  +6.25%  "

$ws.Range("E40").Value = "  -1.89%  "

$ws.Range("E41").Value = "  +2.09%  "

$ws.Range("E42").Value = "  -2.92%  "

$ws.Range("E43").Value = "  -2.11%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "90.72"
$ws.Range("E44").Value = "  -3.05%  "

$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "15.84"
$ws.Range("E45").Value = "  -0.13%  "

$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0880"
$ws.Range("E46").Value = "  -4.21%  "

$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").Value = "1.364.31"
$ws.Range("E47").Value = "  -2.38%  "

$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.46"
$ws.Range("E48").Value = "  -0.76%  "

$ws.Range("E49").Value = "  +1.12%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "46.09"
$ws.Range("E50").Value = "  +2.13%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.65"
$ws.Range("E51").Value = "  +11.46%  "
